$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: extend the existing "Caracteristicas" / Contagem text (column F)
#     with an extra sentence, then add "Principal Assunto" (column E). The
#     F3 update must be applied before the new E3 value so the shared-string
#     table keeps the (modified) F3 text at its original index and appends
#     the brand-new E3 text right after it. ---

$existingF3 = 'O autor questiona como seria possível analisar tendencias emergentes através das estratégias de determinados atores acerca do futuro. Assim, ele prossegue questionando, como seria possível que sinais antecipativos de novas tendências apareçam em um campo que não está formalmente estabilizado? A hipótese colocada pelo autor é que os sinais antecipativos de mudanças futuras não são encontrados em sua forma integra (em inglês ele coloca "are not given once and for all"), mas existem internamente em processos constantes de mudança e poderiam ser mapeados através de mapas de folksonomia (citar o que é). Através de uma pesquisa feita em sites que possuem tags de caracterização para algumas start-ups, os autores fizeram um levantamento de evolução desses termos (tags) ao longo de 2005 e 2007 com fins de verificar como os mesmos termos evoluem do passado (2005) para o futuro (2007), verificando quais eram os seus significados e como eles se desenvolveram. '
$ws.Range("F3").Value = $existingF3

$ws.Range("E3").Value = 'Verificar como as tags são evoluídas e como seus significados se modificam ao longo do tempo.'

$ws.Rows.Item(3).RowHeight = 255

# --- Row 4: brand-new reference (Liebl & Schwarz) ---

$bibtex4 = @'
@article{article,
author = {Liebl, Franz and Schwarz, Jan Oliver},
year = {2010},
month = {05},
pages = {313-327},
title = {Normality of the future: Trend diagnosis for strategic foresight},
volume = {42},
journal = {Futures},
doi = {10.1016/j.futures.2009.11.017}
}
'@

$ws.Range("A4").Value = $bibtex4

$authors4 = @'
Jan Oliver Schwarz
Franz Liebl
'@

$ws.Range("B4").Value = $authors4

$ws.Range("C4").Value = 'Normality of the future: Trend diagnosis for strategic foresight'

$ws.Rows.Item(4).RowHeight = 150

# --- View state: move selection to the new entry row and scroll so row 3
#     is at the top of the visible pane ---

$ws.Range("C4").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
